$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Font.Bold = 0
